$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Add a new row at the end of the table; Word inherits formatting
# (including the Times New Roman rPr used by columns 3 & 4) from the
# preceding row automatically.
$newRow = $t.Rows.Add()
$rowIndex = $t.Rows.Count

# Column 1: plain text
$cell1 = $t.Cell($rowIndex, 1)
$cell1.Range.Text = "Journals with p-values"

# Column 2: hyperlink followed by a trailing space run
$cell2 = $t.Cell($rowIndex, 2)
$rng2 = $cell2.Range
$startPos = $rng2.Start
$url = "https://youtu.be/aUNeLox7920"
$rng2.Text = $url + " "

$linkRange = $d.Range($startPos, $startPos + $url.Length)
$d.Hyperlinks.Add($linkRange, $url, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $url)

# Column 3: "Created" date text (keeps inherited Times New Roman formatting)
$cell3 = $t.Cell($rowIndex, 3)
$cell3.Range.Text = "Created 4/4/24"

# Column 4: Stone reference text (keeps inherited Times New Roman formatting)
$cell4 = $t.Cell($rowIndex, 4)
$cell4.Range.Text = "Stone 52"
